$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before column DR ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Insert a full column before DR - this shifts DR..EV (and their formatting) to ES..EW
$wsPrix.Columns("DR").Insert()

# New column header (row 1) - new daily price column for 13-nov
$wsPrix.Range("DR1").Value = "13-nov"

# Fill data rows (2-25) of the new column with the "-" placeholder used
# throughout the sheet for missing values
$wsPrix.Range("DR2:DR25").Value = "-"

# --- Sheet "Gaz": append new row with latest price ---
# The Date column stores plain text dates (e.g. "2025-11-10"), not real
# Excel dates, so force text formatting before assigning the value to
# avoid automatic date-serial conversion, then restore the neighbouring
# cell's (unstyled) look so the new row matches the existing ones.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A150").NumberFormat = "@"
$wsGaz.Range("A150").Value = "2025-11-11"
$wsGaz.Range("A150").Style = $wsGaz.Range("A149").Style
$wsGaz.Range("B150").Value = 28.7

# --- Sheet "CO2": append new row with latest price ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A150").NumberFormat = "@"
$wsCo2.Range("A150").Value = "2025-11-11"
$wsCo2.Range("A150").Style = $wsCo2.Range("A149").Style
$wsCo2.Range("B150").Value = 80.42
